$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $text) {
    # Force the value to be stored as text (not auto-converted to a number/date),
    # matching Excel's literal "quote-prefix" behavior, then reset the style back
    # to Normal so no stray formatting is left behind on the cell.
    $sheet.Range($addr).Value = "'" + $text
    $sheet.Range($addr).Style = "Normal"
}

Set-TextCell $ws "D2" '19.938.27'
Set-TextCell $ws "E2" '  -8.24%  '
Set-TextCell $ws "D3" '1.402.33'
Set-TextCell $ws "E3" '  -8.94%  '
Set-TextCell $ws "D4" '1.002'
Set-TextCell $ws "E4" '  +0.35%  '
Set-TextCell $ws "D5" '1.003'
Set-TextCell $ws "E5" '  +0.35%  '
Set-TextCell $ws "D6" '273.86'
Set-TextCell $ws "E6" '  -5.63%  '
Set-TextCell $ws "D7" '0.3672'
Set-TextCell $ws "E7" '  -6.99%  '
Set-TextCell $ws "D8" '0.3115'
Set-TextCell $ws "E8" '  -2.89%  '
Set-TextCell $ws "D9" '39.60'
Set-TextCell $ws "E9" '  -7.89%  '
Set-TextCell $ws "D10" '1.008'
Set-TextCell $ws "E10" '  -7.20%  '
Set-TextCell $ws "D11" '0.06480'
Set-TextCell $ws "E11" '  -10.11%  '
Set-TextCell $ws "D12" '1.002'
Set-TextCell $ws "E12" '  +0.33%  '
Set-TextCell $ws "D13" '5.437'
Set-TextCell $ws "E13" '  -5.74%  '
Set-TextCell $ws "D14" '17.33'
Set-TextCell $ws "E14" '  -6.10%  '
Set-TextCell $ws "D15" '6.128'
Set-TextCell $ws "E15" '  -8.01%  '
Set-TextCell $ws "D16" '1.404.95'
Set-TextCell $ws "E16" '  -8.97%  '
Set-TextCell $ws "D17" '0.00001010'
Set-TextCell $ws "E17" '  -8.35%  '
Set-TextCell $ws "D18" '0.05687'
Set-TextCell $ws "E18" '  -13.81%  '
Set-TextCell $ws "D20" '70.42'
Set-TextCell $ws "E20" '  -16.65%  '
Set-TextCell $ws "D21" '5.544'
Set-TextCell $ws "E21" '  -10.09%  '
Set-TextCell $ws "D22" '14.69'
Set-TextCell $ws "E22" '  -6.00%  '
Set-TextCell $ws "D23" '10.99'
Set-TextCell $ws "E23" '  +0.88%  '
Set-TextCell $ws "D24" '2.269'
Set-TextCell $ws "E24" '  -4.31%  '
Set-TextCell $ws "D25" '19.960.08'
Set-TextCell $ws "E25" '  -8.14%  '
Set-TextCell $ws "D26" '2.220'
Set-TextCell $ws "E26" '  -7.99%  '
Set-TextCell $ws "D27" '135.20'
Set-TextCell $ws "E27" '  -11.05%  '
Set-TextCell $ws "D28" '16.84'
Set-TextCell $ws "E28" '  -9.18%  '
Set-TextCell $ws "D29" '1.563.38'
Set-TextCell $ws "E29" '  -8.72%  '
Set-TextCell $ws "D30" '108.76'
Set-TextCell $ws "E30" '  -7.58%  '
Set-TextCell $ws "D31" '4.076'
Set-TextCell $ws "E31" '  -16.03%  '
Set-TextCell $ws "D32" '5.269'
Set-TextCell $ws "E32" '  -14.45%  '
Set-TextCell $ws "D33" '0.8119'
Set-TextCell $ws "E33" '  -17.59%  '
Set-TextCell $ws "D34" '0.07648'
Set-TextCell $ws "E34" '  -6.23%  '
Set-TextCell $ws "D35" '8.381'
Set-TextCell $ws "E35" '  -3.03%  '
Set-TextCell $ws "D36" '1.444'
Set-TextCell $ws "E36" '  -3.00%  '
Set-TextCell $ws "D37" '0.05767'
Set-TextCell $ws "E37" '  -4.25%  '
Set-TextCell $ws "D38" '4.808'
Set-TextCell $ws "E38" '  -8.04%  '
Set-TextCell $ws "D39" '1.001'
Set-TextCell $ws "E39" '  +0.26%  '
Set-TextCell $ws "D40" '0.02065'
Set-TextCell $ws "E40" '  -8.34%  '
Set-TextCell $ws "D41" '0.1898'
Set-TextCell $ws "E41" '  -7.54%  '
Set-TextCell $ws "D42" '10.32'
Set-TextCell $ws "E42" '  -9.37%  '
Set-TextCell $ws "D43" '1.094'
Set-TextCell $ws "E43" '  -8.08%  '
Set-TextCell $ws "D46" '3.513'
Set-TextCell $ws "E46" '  -5.91%  '
Set-TextCell $ws "D47" '0.5098'
Set-TextCell $ws "E47" '  -9.10%  '
Set-TextCell $ws "D48" '111.49'
Set-TextCell $ws "E48" '  -4.79%  '
Set-TextCell $ws "D49" '1.752'
Set-TextCell $ws "E49" '  -8.00%  '
Set-TextCell $ws "E19" '  +0.29%  '

# Rows 44/45 and 50/51 had their coin entries swapped, with updated price/volume
Set-TextCell $ws "B44" 'TheSandbox'
Set-TextCell $ws "C44" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell $ws "D44" '0.5258'
Set-TextCell $ws "E44" '  -10.24%  '
Set-TextCell $ws "B45" 'EnergySwap'
Set-TextCell $ws "C45" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws "D45" '12.24'
Set-TextCell $ws "E45" '  -7.83%  '
Set-TextCell $ws "B50" 'PaxDollar'
Set-TextCell $ws "C50" 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell $ws "D50" '1.001'
Set-TextCell $ws "E50" '  +0.28%  '
Set-TextCell $ws "B51" 'EOS'
Set-TextCell $ws "C51" 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextCell $ws "D51" '1.029'
Set-TextCell $ws "E51" '  -12.15%  '
